$wb = $excel.ActiveWorkbook

# Copy the "Norway" sheet to the end of the workbook to create the "Italy" sheet,
# then copy "Italy" to the end again to create the "Spain" sheet.
$norway = $wb.Worksheets.Item("Norway")
$norway.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"

$italy.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"

# Fill in the market name cells first (both sheets), then the code cells (both sheets),
# matching the order new shared strings were appended in the saved file.
$italy.Range("B2").Value = "Italy Market"
$spain.Range("B2").Value = "Spain Market"

$spain.Range("B4").Value = "NGC-3442/T2128/T2127/T2130"
$italy.Range("B4").Value = "NGC-3443/T1967/T1968/T1970"

# Restore the selections/active cells on each sheet
$spain.Range("A1:XFD1048576").Select() | Out-Null
$italy.Range("A8").Select() | Out-Null
